$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected (all cells locked) - unprotect so the
# refreshed figures below can be written, then restore protection.
$ws.Unprotect()

# Update the confidentiality / "as of" date notice text (A10):
# 2021-05-06 -> 2021-05-07
$newText = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + "`n" + "Model holdings provided as of 2021-05-07 for illustrative purposes only and are subject to change."
$ws.Range("A10").Value = $newText

# Update the recalculated weight / percent-change figures
$ws.Range("D2").Value = 0.4780047497181198
$ws.Range("E2").Value = 0.005080109417741285

$ws.Range("D3").Value = 0.3413382161828364
$ws.Range("E3").Value = 0.005616399887672152

$ws.Range("D4").Value = 0.09576119916557345
$ws.Range("E4").Value = 0.01329500538986705

$ws.Range("D5").Value = 0.05362327300433496
$ws.Range("E5").Value = 0.0006875214850463607

$ws.Range("D6").Value = 0.03127256192913518
$ws.Range("E6").Value = 0.01296787503684049

$ws.Range("E7").Value = 0.006060959836309721

# Restore the sheet protection that was in place before the edit.
$ws.Protect()
